$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "D7"  = -7.453
    "A10" = -21.846
    "A12" = -21.589
    "D15" = -8.175000000000001
    "A18" = -22.403
    "D20" = -7.561999999999999
    "D29" = -7.292
    "D30" = -7.178999999999999
    "D31" = -8.023
    "A37" = -20.21700000000001
    "D40" = -7.663000000000001
    "A55" = -22.283
    "A68" = -21.534
    "D68" = -6.778
    "D76" = -7.753000000000002
    "A77" = -20.934
    "A78" = -20.281
    "D87" = -8.113
    "D88" = -8.030000000000001
    "D96" = -7.267
    "D98" = -8.228
    "D101" = -7.987
    "D102" = -8.142999999999999
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
